$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3-6 are removed from the log, leaving only the header row (1) and one data row (2)
$ws.Rows("3:6").Delete()

# Update the remaining data row (row 2) with the new values
# Force text storage so the numeric-looking phone number stays a string, matching
# the original inline-string cell type.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "201032822563"

$newMessage = @"
Message: no such element: Unable to locate element: {"method":"xpath","selector":"//div[@contenteditable="true"][@data-tab="10"]"}
  (Session info: chrome=141.0.7390.108); For documentation on this error, please visit: https://www.selenium.dev/documentation/webdriver/troubleshooting/errors#no-such-element-exception
Stacktrace:
`tGetHandleVerifier [0x0x7ff607dce9e5+80021]
`tGetHandleVerifier [0x0x7ff607dcea40+80112]
`t(No symbol) [0x0x7ff607b5060f]
`t(No symbol) [0x0x7ff607ba8854]
`t(No symbol) [0x0x7ff607ba8b1c]
`t(No symbol) [0x0x7ff607bfc927]
`t(No symbol) [0x0x7ff607bd126f]
`t(No symbol) [0x0x7ff607bf968a]
`t(No symbol) [0x0x7ff607bd1003]
`t(No symbol) [0x0x7ff607b995d1]
`t(No symbol) [0x0x7ff607b9a3f3]
`tGetHandleVerifier [0x0x7ff60808dd8d+2960445]
`tGetHandleVerifier [0x0x7ff60808804a+2936570]
`tGetHandleVerifier [0x0x7ff6080a8a87+3070263]
`tGetHandleVerifier [0x0x7ff607de84ce+185214]
`tGetHandleVerifier [0x0x7ff607deff1f+216527]
`tGetHandleVerifier [0x0x7ff607dd7c24+117460]
`tGetHandleVerifier [0x0x7ff607dd7ddf+117903]
`tGetHandleVerifier [0x0x7ff607dbdcb8+11112]
`tBaseThreadInitThunk [0x0x7ff9476be8d7+23]
`tRtlUserThreadStart [0x0x7ff94821c34c+44]

"@

$ws.Range("B2").Value = $newMessage

# Setting a multi-line value auto-expands the row height; restore the
# row to its default (un-customized) height so only the cell content changed.
$ws.Rows(2).AutoFit()
